# Auto-generated edit script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is purely numeric-looking (e.g. "1.00", "546.49")
# must be forced to Text format first, otherwise Excel COM auto-converts the
# assigned string into a real number and the formatting (trailing zeros, etc.) is lost.
$forceTextCells = @("D4", "D5", "D6", "D7", "D14", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D35", "D36", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value changes described by the diff
$ws.Range("D2").Value = "60.948.34"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "2.369.20"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "546.49"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "132.52"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").Value = "2.369.97"
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "24.13"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "2.791.67"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("D16").Value = "60.883.64"
$ws.Range("E16").Value = "  +5.41%  "
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "2.340.24"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "7.04"
$ws.Range("E20").Value = "  +10.42%  "
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "317.86"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "63.47"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +4.51%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "8.06"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").Value = "1.36"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "171.78"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.73"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("E32").Value = "  +8.66%  "
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("E34").Value = "  +14.39%  "
$ws.Range("D35").Value = "0.384"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Value = "18.14"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +9.02%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "319.34"
$ws.Range("E40").Value = "  +10.65%  "
$ws.Range("D41").Value = "1.55"
$ws.Range("E41").Value = "  +4.85%  "
$ws.Range("D42").Value = "38.39"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "143.85"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "3.48"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "19.54"
$ws.Range("E46").Value = "  +8.39%  "
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "0.566"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "0.0₆0206"
$ws.Range("E51").Value = "  +1.90%  "
